$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for several rows to repulled/recomputed data
$ws.Range("F3").Value = -5
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = -7
$ws.Range("F14").Value = 2
